# record counter and finalDesign
#
# For each of the three worksheets, move the "list label" (e.g.
# WORDS_WITH_PHTYPE_PER_LIST) from C4 up into the merged header cell C1,
# then add two new summary rows ("Big" / "Small") of VOWEL statistics
# below the existing row, matching the layout/format already used by
# the sheet (percentage columns D:K, with a right border closing the
# "Height" (D:H) and "Backness" (I:K) sub-groups).

$wb = $excel.ActiveWorkbook

function Fill-SummaryRow($ws, $rowNum, $label, $values, $numFmt) {
    $ws.Range("A" + $rowNum + ":C" + $rowNum).ClearFormats()
    $ws.Range("A" + $rowNum).Value = $label
    $ws.Range("B" + $rowNum).Value = 30
    $ws.Range("C" + $rowNum).Value = ""

    $ws.Range("A" + $rowNum).Borders.Item(10).LineStyle = 1
    $ws.Range("B" + $rowNum).Borders.Item(10).LineStyle = 1
    $ws.Range("C" + $rowNum).Borders.Item(10).LineStyle = 1

    $cols = @("D", "E", "F", "G", "H", "I", "J", "K")
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $cell = $ws.Range($cols[$i] + $rowNum)
        $cell.Value = $values[$i]
        $cell.NumberFormat = $numFmt
    }

    $ws.Range("H" + $rowNum).Borders.Item(10).LineStyle = 1
    $ws.Range("K" + $rowNum).Borders.Item(10).LineStyle = 1
}

function Process-Sheet($ws, $listLabel, $row4Values, $row5Values, $numFmt) {
    # Move the list-name label from C4 up into the merged C1 header cell.
    $ws.Range("C1").Value = $listLabel

    # The existing data row becomes the "Big" summary row; its label used
    # to live in C4, now C4 is cleared out (blank, matching A4/B4).
    $ws.Range("A4:C4").ClearFormats()
    $ws.Range("A4").Value = "Big"
    $ws.Range("B4").Value = 30
    $ws.Range("C4").Value = ""

    $ws.Range("A4").Borders.Item(10).LineStyle = 1
    $ws.Range("B4").Borders.Item(10).LineStyle = 1
    $ws.Range("C4").Borders.Item(10).LineStyle = 1

    $cols = @("D", "E", "F", "G", "H", "I", "J", "K")
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $cell = $ws.Range($cols[$i] + "4")
        $cell.Value = $row4Values[$i]
    }

    $ws.Range("H4").Borders.Item(10).LineStyle = 1
    $ws.Range("K4").Borders.Item(10).LineStyle = 1

    # New "Small" summary row.
    Fill-SummaryRow $ws 5 "Small" $row5Values $numFmt
}

$ws1 = $wb.Worksheets.Item(1)
Process-Sheet $ws1 "WORDS_WITH_PHTYPE_PER_LIST" `
    @(0.4, 0.03333333333333333, 0.03333333333333333, 0.5, 0.4, 0.7666666666666667, 0.06666666666666667, 0.43333333333333335) `
    @(0.3333333333333333, 0.06666666666666667, 0.0, 0.36666666666666664, 0.7333333333333333, 0.8333333333333334, 0.0, 0.36666666666666664) `
    "0.0%"

$ws2 = $wb.Worksheets.Item(2)
Process-Sheet $ws2 "PHTYPES_PER_LIST" `
    @(0.11333333333333333, 0.006666666666666667, 0.013333333333333334, 0.16, 0.13333333333333333, 0.26, 0.02666666666666667, 0.14) `
    @(0.0684931506849315, 0.0136986301369863, 0.0, 0.0958904109589041, 0.19863013698630136, 0.2876712328767123, 0.0, 0.08904109589041095) `
    "0.0%"

$ws3 = $wb.Worksheets.Item(3)
Process-Sheet $ws3 "PHTYPES_AVERAGE_PER_WORD" `
    @(0.5666666666666667, 0.03333333333333333, 0.06666666666666667, 0.8, 0.6666666666666666, 1.3, 0.13333333333333333, 0.7) `
    @(0.3333333333333333, 0.06666666666666667, 0.0, 0.4666666666666667, 0.9666666666666667, 1.4, 0.0, 0.43333333333333335) `
    "0.0"
